# Update "想去人数" (wanted-to-go attendance count) values in the "展览"
# and "全部类型" worksheets to match the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")

$sheet1Updates = @{
    "F2"  = 608
    "F3"  = 45
    "F5"  = 14
    "F6"  = 15115
    "F7"  = 408
    "F9"  = 683
    "F10" = 15265
    "F11" = 43
    "F12" = 8795
    "F13" = 345
    "F14" = 5
    "F15" = 72
    "F18" = 185
    "F19" = 16
    "F20" = 31
    "F21" = 524
    "F23" = 7
    "F24" = 53
    "F25" = 1090
    "F26" = 10
    "F27" = 14
    "F28" = 59
    "F29" = 29
    "F31" = 415
    "F32" = 33
    "F33" = 31
    "F34" = 233
    "F35" = 280
    "F36" = 431
    "F38" = 5404
    "F39" = 5227
}

foreach ($addr in $sheet1Updates.Keys) {
    $ws1.Range($addr).Value = $sheet1Updates[$addr]
}

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4Updates = @{
    "F2"  = 608
    "F3"  = 45
    "F5"  = 14
    "F6"  = 15115
    "F7"  = 408
    "F9"  = 683
    "F10" = 15265
    "F11" = 43
    "F12" = 8795
    "F13" = 345
    "F14" = 5
    "F16" = 72
    "F19" = 185
    "F20" = 16
    "F21" = 31
    "F22" = 524
    "F24" = 7
    "F25" = 53
    "F26" = 1090
    "F27" = 10
    "F28" = 14
    "F29" = 59
    "F30" = 29
    "F34" = 415
    "F35" = 33
    "F36" = 31
    "F37" = 233
    "F38" = 280
    "F39" = 431
    "F41" = 5404
    "F42" = 5227
}

foreach ($addr in $sheet4Updates.Keys) {
    $ws4.Range($addr).Value = $sheet4Updates[$addr]
}
